$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date in column C for all data rows (2-497)
#    from 45189 to 45190.
for ($r = 2; $r -le 497; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}

# 2. Row 497 picks up an explicit row height stamp.
$ws.Rows.Item(497).RowHeight = 15

# 3. Append the new record as row 498.
$ws.Cells.Item(498, 1).Value = "A 44001-2023"
$ws.Cells.Item(498, 2).Value = 45188
$ws.Cells.Item(498, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(498, 3).Value = 45190
$ws.Cells.Item(498, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(498, 4).Value = "KALMAR LÄN"
$ws.Cells.Item(498, 5).Value = "MÖNSTERÅS"
$ws.Cells.Item(498, 7).Value = 1.9
$ws.Cells.Item(498, 8).Value = 0
$ws.Cells.Item(498, 9).Value = 0
$ws.Cells.Item(498, 10).Value = 0
$ws.Cells.Item(498, 11).Value = 0
$ws.Cells.Item(498, 12).Value = 0
$ws.Cells.Item(498, 13).Value = 0
$ws.Cells.Item(498, 14).Value = 0
$ws.Cells.Item(498, 15).Value = 0
$ws.Cells.Item(498, 16).Value = 0
$ws.Cells.Item(498, 17).Value = 0
$ws.Cells.Item(498, 18).Value = ""
$ws.Cells.Item(498, 18).WrapText = $true
